# Make room for logos: hide row 4, mark customHeight on several rows,
# and move the active selection from E10 to F10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert explicit row heights (causes customHeight="1" to be written)
# while keeping the existing heights unchanged. Hidden rows report a
# RowHeight of 0 via COM, so use the known original heights instead of
# round-tripping the current value.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 105
$ws.Rows.Item(7).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 45

# Row 4 becomes newly hidden (to make room for logos).
$ws.Rows.Item(4).Hidden = $true

# Move the selection from E10 to F10.
$ws.Range("F10").Select()
